# This script updates the species/woodiness mapping sheet so that:
#  - family-level rows (e.g. "Apiaceae") always keep a "<Family> species" label,
#  - genus-level rows that previously held only a bare genus name or a
#    GBIF-assigned species name are updated so that the *original* species
#    names supplied by the user are preserved (shifted into place) and a
#    trailing "<Genus> species" row is kept for whichever mapping could not
#    be tied to one specific species.
# The concrete before/after values below were derived from the authoritative
# diff of the workbook's XML and are applied here as direct cell writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "Acer campestre"
$ws.Cells.Item(3, 2).Value = "woody"
$ws.Cells.Item(4, 1).Value = "Acer negundo"
$ws.Cells.Item(5, 1).Value = "Acer platanoides"
$ws.Cells.Item(6, 1).Value = "Acer pseudoplatanus"
$ws.Cells.Item(7, 1).Value = "Acer species"
$ws.Cells.Item(7, 2).Value = "not found"

$ws.Cells.Item(11, 1).Value = "Allium species"
$ws.Cells.Item(11, 2).Value = "not found"

$ws.Cells.Item(17, 1).Value = "Apiaceae species"

$ws.Cells.Item(22, 1).Value = "Asteraceae species"
$ws.Cells.Item(22, 2).Value = "not found"

$ws.Cells.Item(27, 1).Value = "Betula species"

$ws.Cells.Item(32, 1).Value = "Brassicaceae species"
$ws.Cells.Item(32, 2).Value = "not found"

$ws.Cells.Item(66, 1).Value = "Draba species"
$ws.Cells.Item(66, 2).Value = "not found"

$ws.Cells.Item(77, 1).Value = "Festuca ovina"
$ws.Cells.Item(78, 1).Value = "Festuca rubra"
$ws.Cells.Item(79, 1).Value = "Festuca rupicola"
$ws.Cells.Item(80, 1).Value = "Festuca species"
$ws.Cells.Item(80, 2).Value = "not found"

$ws.Cells.Item(91, 1).Value = "Geranium molle"
$ws.Cells.Item(91, 2).Value = "herbaceous"
$ws.Cells.Item(92, 1).Value = "Geranium pratense"
$ws.Cells.Item(93, 1).Value = "Geranium pusillum"
$ws.Cells.Item(94, 1).Value = "Geranium pyrenaicum"
$ws.Cells.Item(94, 2).Value = "not found"
$ws.Cells.Item(95, 1).Value = "Geranium rotundifolium"
$ws.Cells.Item(95, 2).Value = "herbaceous"
$ws.Cells.Item(96, 1).Value = "Geranium species"
$ws.Cells.Item(96, 2).Value = "not found"

$ws.Cells.Item(123, 1).Value = "Medicago falcata"
$ws.Cells.Item(123, 2).Value = "herbaceous"
$ws.Cells.Item(124, 1).Value = "Medicago lupulina"
$ws.Cells.Item(125, 1).Value = "Medicago species"
$ws.Cells.Item(125, 2).Value = "not found"

$ws.Cells.Item(145, 1).Value = "Poaceae species"
$ws.Cells.Item(145, 2).Value = "not found"

$ws.Cells.Item(151, 1).Value = "Prunus avium"
$ws.Cells.Item(152, 1).Value = "Prunus mahaleb"
$ws.Cells.Item(153, 1).Value = "Prunus species"
$ws.Cells.Item(153, 2).Value = "not found"

$ws.Cells.Item(159, 1).Value = "Rubus caesius"
$ws.Cells.Item(159, 2).Value = "woody"
$ws.Cells.Item(160, 1).Value = "Rubus idaeus"
$ws.Cells.Item(161, 1).Value = "Rubus species"
$ws.Cells.Item(161, 2).Value = "not found"

$ws.Cells.Item(167, 1).Value = "Senecio jacobaea"
$ws.Cells.Item(167, 2).Value = "not found"
$ws.Cells.Item(168, 1).Value = "Senecio species"

$ws.Cells.Item(202, 1).Value = "Triticum species"
$ws.Cells.Item(202, 2).Value = "not found"

$ws.Cells.Item(213, 1).Value = "Vicia cracca"
$ws.Cells.Item(213, 2).Value = "herbaceous"
$ws.Cells.Item(214, 1).Value = "Vicia hirsuta"
$ws.Cells.Item(215, 1).Value = "Vicia sativa"
$ws.Cells.Item(216, 1).Value = "Vicia sepium"
$ws.Cells.Item(216, 2).Value = "not assigned"
$ws.Cells.Item(217, 1).Value = "Vicia species"
$ws.Cells.Item(217, 2).Value = "not found"
